# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.256.58"
$ws.Range("E2").Value = "  +5.57%  "
$ws.Range("D3").Value = "3.275.86"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'407.32"
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("D6").Value = "'111.41"
$ws.Range("E6").Value = "  +3.13%  "
$ws.Range("D7").Value = "3.275.94"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D8").Value = "'0.564"
$ws.Range("E8").Value = "  -3.05%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "'0.615"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "'0.113"
$ws.Range("E11").Value = "  +13.75%  "
$ws.Range("D12").Value = "'38.54"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "3.782.19"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "'8.13"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "'18.89"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "3.280.80"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "59.975.14"
$ws.Range("E18").Value = "  +5.33%  "
$ws.Range("D19").Value = "'0.979"
$ws.Range("E19").Value = "  -4.61%  "
$ws.Range("D20").Value = "'10.55"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").Value = "'3.27"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("D23").Value = "'12.40"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("D24").Value = "'295.40"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'73.15"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").Value = "'3.07"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").Value = "'28.84"
$ws.Range("E27").Value = "  +3.27%  "
$ws.Range("E28").Value = "  -2.48%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.39"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.172"
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("D31").Value = "'7.44"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").Value = "'0.113"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "'11.12"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("E35").Value = "  +14.28%  "
$ws.Range("D36").Value = "'38.73"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").Value = "'0.0477"
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("D38").Value = "'52.07"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").Value = "'0.996"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").Value = "  +6.25%  "
$ws.Range("D41").Value = "'3.29"
$ws.Range("E41").Value = "  -4.80%  "
$ws.Range("D42").Value = "'134.79"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "'0.292"
$ws.Range("E43").Value = "  +4.89%  "
$ws.Range("D44").Value = "'0.120"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").Value = "'16.16"
$ws.Range("E46").Value = "  -4.21%  "
$ws.Range("D47").Value = "'3.75"
$ws.Range("E47").Value = "  -4.60%  "
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D49").Value = "'20.83"
$ws.Range("E49").Value = "  -5.05%  "
$ws.Range("D50").Value = "2.109.58"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "3.598.03"
$ws.Range("E51").Value = "  +0.75%  "
